$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new column labels (template for bulk import) ---
$ws.Range("A1").Value = "Nombre (texto)"
$ws.Range("B1").Value = "ID Tipo (Numero"
$ws.Range("C1").Value = "Fecha (AAAA-MM-DD)"
$ws.Range("D1").Value = "Porcentaje Implementacion (numero)"
$ws.Range("E1").Value = "ID Empleado que elaboro (numero)"
$ws.Range("F1").Value = "Estatus (Texto)"

# Recolor the header highlight fill from green to yellow. Set it on A1
# directly, then fan the same format out to the rest of the header with a
# format-only paste so the shared date format on C1 isn't disturbed.
$ws.Range("A1").Interior.Color = 65535
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New empty, underlined cell next to the header (F2)
$ws.Range("F2").Font.Underline = $true

# --- Column widths to fit the new, longer header text ---
$ws.Columns.Item(1).ColumnWidth = 13.333333
$ws.Columns.Item(2).ColumnWidth = 14
$ws.Columns.Item(3).ColumnWidth = 18.833333
$ws.Columns.Item(4).ColumnWidth = 32.5
$ws.Columns.Item(5).ColumnWidth = 30.333333
$ws.Columns.Item(6).ColumnWidth = 12.666667

# --- Selection parks on D21 ---
$ws.Range("D21").Select() | Out-Null

# --- Page setup: portrait, standard paper ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
